$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new BOM line for the standoffs (row 22): Qty, Vendor, Part Number, Price
$ws.Range("D22").Value = 4
$ws.Range("E22").Value = "McMaster"
$ws.Range("F22").Value = "92745A326"
$ws.Range("G22").Value = 0.46

# Extend the Total formula (Price * Qty) down through the new rows 21 and 22
$ws.Range("H21").Formula = "=G21*D21"
$ws.Range("H22").Formula = "=G22*D22"

# Restore the selection that was active when the workbook was saved
$ws.Range("H18").Select() | Out-Null
